$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()),
                (''selector'',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'',
                                                     random_state=42))),
                (''model'',
                 AdaBoostClassifier(estimator=RandomForestClassifier(max_depth=1,
                                                                     min_samples_leaf=2,
                                                                     n_estimators=10,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B2").Value = 0.6356809856809857
$ws.Range("C2").Value = '{''selector'': SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'', random_state=42)), ''scaler'': StandardScaler(), ''model__n_estimators'': 5, ''model__estimator__n_estimators'': 10, ''model__estimator__min_samples_split'': 2, ''model__estimator__min_samples_leaf'': 2, ''model__estimator__max_features'': ''sqrt'', ''model__estimator__max_depth'': 1, ''model__estimator__class_weight'': None}'
$ws.Range("D2").Value = 0.9650100618908158
$ws.Range("E2").Value = 0.5596665279165279
$ws.Range("F2").Value = 0.7096774193548386
$ws.Range("G2").Value = 0.9641279004658485
$ws.Range("H2").Value = 0.5549948412698412
$ws.Range("I2").Value = 0.7333333333333333
$ws.Range("J2").Value = 0.966872340425532
$ws.Range("K2").Value = 0.5863333333333334
$ws.Range("L2").Value = 0.6875
$ws.Range("M2").Value = '[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]'
$ws.Range("N2").Value = '[1 0 1 1 1 1 1 0 1 0 1 1 0 1 0 0 0 1 0 1 1 1 1 0]'

$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faaa3d730d0>),
                (''model'',
                 AdaBoostClassifier(estimator=RandomForestClassifier(max_depth=3,
                                                                     max_features=''log2'',
                                                                     min_samples_leaf=5,
                                                                     min_samples_split=5,
                                                                     n_estimators=50,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B3").Value = 0.6490276390276389
$ws.Range("C3").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa0022e040>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 5, ''model__estimator__n_estimators'': 50, ''model__estimator__min_samples_split'': 5, ''model__estimator__min_samples_leaf'': 5, ''model__estimator__max_features'': ''log2'', ''model__estimator__max_depth'': 3, ''model__estimator__class_weight'': None}'
$ws.Range("D3").Value = 0.9772803325874193
$ws.Range("E3").Value = 0.5341145743145743
$ws.Range("F3").Value = 0.7222222222222223
$ws.Range("G3").Value = 0.9757917630066838
$ws.Range("H3").Value = 0.5882380952380953
$ws.Range("I3").Value = 0.65
$ws.Range("J3").Value = 0.9793191489361702
$ws.Range("K3").Value = 0.5111666666666668
$ws.Range("L3").Value = 0.8125
$ws.Range("M3").Value = '[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]'
$ws.Range("N3").Value = '[1 1 1 1 1 1 0 1 1 1 1 0 0 0 1 1 1 1 1 1 1 1 1 1]'

$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f7e6460>),
                (''model'',
                 AdaBoostClassifier(estimator=RandomForestClassifier(max_depth=2,
                                                                     max_features=''log2'',
                                                                     min_samples_leaf=4,
                                                                     min_samples_split=5,
                                                                     n_estimators=50,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B4").Value = 0.6963347763347764
$ws.Range("C4").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa001ca4c0>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 5, ''model__estimator__n_estimators'': 50, ''model__estimator__min_samples_split'': 5, ''model__estimator__min_samples_leaf'': 4, ''model__estimator__max_features'': ''log2'', ''model__estimator__max_depth'': 2, ''model__estimator__class_weight'': None}'
$ws.Range("D4").Value = 0.9768167183757286
$ws.Range("E4").Value = 0.5672153513153512
$ws.Range("F4").Value = 0.6451612903225806
$ws.Range("G4").Value = 0.9761439665258788
$ws.Range("H4").Value = 0.5861103174603175
$ws.Range("I4").Value = 0.8333333333333334
$ws.Range("J4").Value = 0.9779111111111112
$ws.Range("K4").Value = 0.5720000000000001
$ws.Range("L4").Value = 0.5263157894736842
$ws.Range("M4").Value = '[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]'
$ws.Range("N4").Value = '[0 1 1 0 0 1 0 1 1 0 0 0 0 1 0 1 0 1 1 0 1 0 1 1]'

$ws.Range("A5").Value = 'Pipeline(steps=[(''scaler'', MinMaxScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f403190>),
                (''model'',
                 AdaBoostClassifier(estimator=RandomForestClassifier(max_depth=4,
                                                                     min_samples_leaf=5,
                                                                     min_samples_split=5,
                                                                     n_estimators=50,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B5").Value = 0.71003996003996
$ws.Range("C5").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa00194130>, ''scaler'': MinMaxScaler(), ''model__n_estimators'': 5, ''model__estimator__n_estimators'': 50, ''model__estimator__min_samples_split'': 5, ''model__estimator__min_samples_leaf'': 5, ''model__estimator__max_features'': ''sqrt'', ''model__estimator__max_depth'': 4, ''model__estimator__class_weight'': None}'
$ws.Range("D5").Value = 0.9811003121169937
$ws.Range("E5").Value = 0.6072379342879342
$ws.Range("F5").Value = 0.7142857142857143
$ws.Range("G5").Value = 0.9771391202975537
$ws.Range("H5").Value = 0.5746496031746031
$ws.Range("I5").Value = 0.7142857142857143
$ws.Range("J5").Value = 0.9857755102040816
$ws.Range("K5").Value = 0.6613333333333332
$ws.Range("L5").Value = 0.7142857142857143
$ws.Range("M5").Value = '[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]'
$ws.Range("N5").Value = '[0 1 1 1 1 1 0 0 0 0 1 1 0 0 1 0 0 1 1 1 0 1 1 1]'

$ws.Range("A6").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 <__main__.NamedFeatureSelector object at 0x7faa9f73c0d0>),
                (''model'',
                 AdaBoostClassifier(estimator=RandomForestClassifier(max_depth=1,
                                                                     max_features=''log2'',
                                                                     min_samples_split=4,
                                                                     n_estimators=10,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B6").Value = 0.749069264069264
$ws.Range("C6").Value = '{''selector'': <__main__.NamedFeatureSelector object at 0x7faa0022e940>, ''scaler'': RobustScaler(), ''model__n_estimators'': 5, ''model__estimator__n_estimators'': 10, ''model__estimator__min_samples_split'': 4, ''model__estimator__min_samples_leaf'': 1, ''model__estimator__max_features'': ''log2'', ''model__estimator__max_depth'': 1, ''model__estimator__class_weight'': None}'
$ws.Range("D6").Value = 0.9726150643258477
$ws.Range("E6").Value = 0.6447921356421357
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.9715124453198313
$ws.Range("H6").Value = 0.6158626984126984
$ws.Range("I6").Value = 0.5
$ws.Range("J6").Value = 0.9746923076923076
$ws.Range("K6").Value = 0.6975
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = '[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]'
$ws.Range("N6").Value = '[1 1 1 1 1 1 1 1 1 0 1 1 1 1 1 0 1 1 1 1 1 1 1 1]'
